# Fipe.xlsx — add the Marca/Modelo/Ano lookup table and format its header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "MarcaSelecionada"
$ws.Cells.Item(1,2).Value = "ModeloSelecionado"
$ws.Cells.Item(1,3).Value = "AnoSelecionado"

# --- Data, written one column at a time (Marca, then Modelo, then Ano) so
# the shared-string table is populated in the same order the source sheet
# was originally authored in (whole columns pasted in sequence).
$colA = @("Acura","Acura","Acura","Acura","Agrale","Agrale","Agrale","Agrale","Alfa Romeo","Alfa Romeo","Alfa Romeo","Alfa Romeo")
$colB = @("Integra GS 1.8","Integra GS 1.8","Legend 3.2/3.5","Legend 3.2/3.5","MARRUÁ 2.8 12V 132cv TDI Diesel","MARRUÁ 2.8 12V 132cv TDI Diesel","MARRUÁ AM 100 2.8  CS TDI Diesel","MARRUÁ AM 100 2.8  CS TDI Diesel","145 Elegant 1.7/1.8 16V","145 Elegant 1.7/1.8 16V","145 Elegant 2.0 16V","145 Elegant 2.0 16V")
$colC = @("1992 Gasolina","1991 Gasolina","1998 Gasolina","1997 Gasolina","2007 Diesel","2006 Diesel","2007 Diesel","2006 Diesel","1999 Gasolina","1998 Gasolina","1998 Gasolina","1997 Gasolina")

for ($i = 0; $i -lt $colA.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt $colC.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}

# --- Header formatting: bold, thin box border, centered/top-aligned -------
# Build the format once on A1, then fan it out with PasteSpecial so every
# header cell shares a single style record instead of each property-set
# minting its own cellXfs entry.
$hdr = $ws.Cells.Item(1,1)
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Copy() | Out-Null
$ws.Range("B1:C1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column widths (best fit) ----------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 29.333333333333332
$ws.Columns.Item(3).ColumnWidth = 14

# --- Selection, matching the saved UI state --------------------------------
$ws.Range("G8").Select() | Out-Null
